$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 62

$ws.Cells.Item($row, 1).Value = "Kindergarden"
$ws.Cells.Item($row, 2).Value = "Kindergarden Harderwijk Bazuindreef"
$ws.Cells.Item($row, 3).Value = "KDV"
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "2024-09-24"
$ws.Cells.Item($row, 4).ClearFormats()
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
